# Applies the data update described by the diff: rows 2-15 and 17-24 of
# Sheet1 get their Fecha (D), Volumen (J), Precio minimo (K), Precio maximo
# (L), Precio promedio ponderado (M) and Precio $/Kg (P) columns updated to
# reflect a refreshed weekly data pull (row 16 and the header are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D, J, K, L, M, P)
$rows = @{
    2  = @(44365, 55, 5000, 5000, 5000, 5000)
    3  = @(44649, 20, 5000, 5000, 5000, 5000)
    4  = @(44679, 50, 5000, 5000, 5000, 5000)
    5  = @(44749, 65, 6000, 6000, 6000, 6000)
    6  = @(44176, 10, 4000, 4000, 4000, 4000)
    7  = @(44280, 55, 4000, 4000, 4000, 4000)
    8  = @(44390, 55, 6000, 6000, 6000, 6000)
    9  = @(44497, 20, 4000, 4000, 4000, 4000)
    10 = @(44656, 85, 5000, 5000, 5000, 5000)
    11 = @(44508, 30, 4000, 4000, 4000, 4000)
    12 = @(44291, 35, 4000, 4000, 4000, 4000)
    13 = @(44313, 20, 4000, 4000, 4000, 4000)
    14 = @(44259, 30, 4000, 4000, 4000, 4000)
    15 = @(44312, 50, 4000, 4000, 4000, 4000)
    17 = @(44315, 40, 4000, 4000, 4000, 4000)
    18 = @(44680, 20, 5000, 5000, 5000, 5000)
    19 = @(44777, 25, 5000, 5000, 5000, 5000)
    20 = @(44316, 20, 4000, 4000, 4000, 4000)
    21 = @(44301, 40, 3000, 3000, 3000, 3000)
    22 = @(44498, 40, 4000, 4000, 4000, 4000)
    23 = @(44781, 40, 5000, 5000, 5000, 5000)
    24 = @(44504, 55, 4000, 4000, 4000, 4000)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
